$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.ClearContents()

# Header row (row 1) - rewritten in original order so shared-string indices 0-19 stay identical
$ws.Cells.Item(1,1).Value = "Sending cluster"
$ws.Cells.Item(1,2).Value = "Ligand symbol"
$ws.Cells.Item(1,3).Value = "Receptor symbol"
$ws.Cells.Item(1,4).Value = "Target cluster"
$ws.Cells.Item(1,5).Value = "Ligand-expressing cells"
$ws.Cells.Item(1,6).Value = "Ligand detection rate"
$ws.Cells.Item(1,7).Value = "Ligand average expression value"
$ws.Cells.Item(1,8).Value = "Ligand total expression value"
$ws.Cells.Item(1,9).Value = "Ligand derived specificity of average expression value"
$ws.Cells.Item(1,10).Value = "Ligand derived specificity of total expression value"
$ws.Cells.Item(1,11).Value = "Receptor-expressing cells"
$ws.Cells.Item(1,12).Value = "Receptor detection rate"
$ws.Cells.Item(1,13).Value = "Receptor average expression value"
$ws.Cells.Item(1,14).Value = "Receptor total expression value"
$ws.Cells.Item(1,15).Value = "Receptor derived specificity of average expression value"
$ws.Cells.Item(1,16).Value = "Receptor derived specificity of total expression value"
$ws.Cells.Item(1,17).Value = "Edge average expression weight"
$ws.Cells.Item(1,18).Value = "Edge total expression weight"
$ws.Cells.Item(1,19).Value = "Edge average expression derived specificity"
$ws.Cells.Item(1,20).Value = "Edge total expression derived specificity"

# Data rows 2-10
# row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Inhbb"
$ws.Cells.Item(2,3).Value = "Acvr1b"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.4201923333333333
$ws.Cells.Item(2,8).Value = 1.260577
$ws.Cells.Item(2,9).Value = 0.08716480679187069
$ws.Cells.Item(2,10).Value = 0.08716480679187069
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 4.453045666666667
$ws.Cells.Item(2,14).Value = 13.359137
$ws.Cells.Item(2,15).Value = 0.4394129038053478
$ws.Cells.Item(2,16).Value = 0.4394129038053478
$ws.Cells.Item(2,17).Value = 1.871135649116556
$ws.Cells.Item(2,18).Value = 16.840220842049
$ws.Cells.Item(2,19).Value = 0.03830134086204801
$ws.Cells.Item(2,20).Value = 0.038301340862048
# row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Inhbb"
$ws.Cells.Item(3,3).Value = "Acvr1b"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.4201923333333333
$ws.Cells.Item(3,8).Value = 1.260577
$ws.Cells.Item(3,9).Value = 0.08716480679187069
$ws.Cells.Item(3,10).Value = 0.08716480679187069
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.150099666666667
$ws.Cells.Item(3,14).Value = 9.450299000000001
$ws.Cells.Item(3,15).Value = 0.3108421842981904
$ws.Cells.Item(3,16).Value = 0.3108421842981904
$ws.Cells.Item(3,17).Value = 1.323647729169222
$ws.Cells.Item(3,18).Value = 11.912829562523
$ws.Cells.Item(3,19).Value = 0.02709449893711483
$ws.Cells.Item(3,20).Value = 0.02709449893711482
# row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Inhbb"
$ws.Cells.Item(4,3).Value = "Acvr1b"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.4201923333333333
$ws.Cells.Item(4,8).Value = 1.260577
$ws.Cells.Item(4,9).Value = 0.08716480679187069
$ws.Cells.Item(4,10).Value = 0.08716480679187069
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.530935
$ws.Cells.Item(4,14).Value = 7.592805
$ws.Cells.Item(4,15).Value = 0.2497449118964618
$ws.Cells.Item(4,16).Value = 0.2497449118964618
$ws.Cells.Item(4,17).Value = 1.063479483165
$ws.Cells.Item(4,18).Value = 9.571315348485001
$ws.Cells.Item(4,19).Value = 0.02176896699270786
$ws.Cells.Item(4,20).Value = 0.02176896699270786
# row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Inhbb"
$ws.Cells.Item(5,3).Value = "Acvr1b"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 2.446732
$ws.Cells.Item(5,8).Value = 7.340196000000001
$ws.Cells.Item(5,9).Value = 0.5075507217365239
$ws.Cells.Item(5,10).Value = 0.5075507217365239
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 4.453045666666667
$ws.Cells.Item(5,14).Value = 13.359137
$ws.Cells.Item(5,15).Value = 0.4394129038053478
$ws.Cells.Item(5,16).Value = 0.4394129038053478
$ws.Cells.Item(5,17).Value = 10.89540933009467
$ws.Cells.Item(5,18).Value = 98.05868397085202
$ws.Cells.Item(5,19).Value = 0.223024336466746
$ws.Cells.Item(5,20).Value = 0.223024336466746
# row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Inhbb"
$ws.Cells.Item(6,3).Value = "Acvr1b"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2.446732
$ws.Cells.Item(6,8).Value = 7.340196000000001
$ws.Cells.Item(6,9).Value = 0.5075507217365239
$ws.Cells.Item(6,10).Value = 0.5075507217365239
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 3.150099666666667
$ws.Cells.Item(6,14).Value = 9.450299000000001
$ws.Cells.Item(6,15).Value = 0.3108421842981904
$ws.Cells.Item(6,16).Value = 0.3108421842981904
$ws.Cells.Item(6,17).Value = 7.707449657622669
$ws.Cells.Item(6,18).Value = 69.36704691860402
$ws.Cells.Item(6,19).Value = 0.1577681749867041
$ws.Cells.Item(6,20).Value = 0.1577681749867041
# row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Inhbb"
$ws.Cells.Item(7,3).Value = "Acvr1b"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.446732
$ws.Cells.Item(7,8).Value = 7.340196000000001
$ws.Cells.Item(7,9).Value = 0.5075507217365239
$ws.Cells.Item(7,10).Value = 0.5075507217365239
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.530935
$ws.Cells.Item(7,14).Value = 7.592805
$ws.Cells.Item(7,15).Value = 0.2497449118964618
$ws.Cells.Item(7,16).Value = 0.2497449118964618
$ws.Cells.Item(7,17).Value = 6.192519654420001
$ws.Cells.Item(7,18).Value = 55.73267688978001
$ws.Cells.Item(7,19).Value = 0.1267582102830738
$ws.Cells.Item(7,20).Value = 0.1267582102830738
# row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Inhbb"
$ws.Cells.Item(8,3).Value = "Acvr1b"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 1.953740666666667
$ws.Cells.Item(8,8).Value = 5.861222
$ws.Cells.Item(8,9).Value = 0.4052844714716054
$ws.Cells.Item(8,10).Value = 0.4052844714716054
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 4.453045666666667
$ws.Cells.Item(8,14).Value = 13.359137
$ws.Cells.Item(8,15).Value = 0.4394129038053478
$ws.Cells.Item(8,16).Value = 0.4394129038053478
$ws.Cells.Item(8,17).Value = 8.700096409490445
$ws.Cells.Item(8,18).Value = 78.300867685414
$ws.Cells.Item(8,19).Value = 0.1780872264765538
$ws.Cells.Item(8,20).Value = 0.1780872264765538
# row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Inhbb"
$ws.Cells.Item(9,3).Value = "Acvr1b"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 1.953740666666667
$ws.Cells.Item(9,8).Value = 5.861222
$ws.Cells.Item(9,9).Value = 0.4052844714716054
$ws.Cells.Item(9,10).Value = 0.4052844714716054
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.150099666666667
$ws.Cells.Item(9,14).Value = 9.450299000000001
$ws.Cells.Item(9,15).Value = 0.3108421842981904
$ws.Cells.Item(9,16).Value = 0.3108421842981904
$ws.Cells.Item(9,17).Value = 6.154477822819779
$ws.Cells.Item(9,18).Value = 55.39030040537801
$ws.Cells.Item(9,19).Value = 0.1259795103743715
$ws.Cells.Item(9,20).Value = 0.1259795103743714
# row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Inhbb"
$ws.Cells.Item(10,3).Value = "Acvr1b"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 1.953740666666667
$ws.Cells.Item(10,8).Value = 5.861222
$ws.Cells.Item(10,9).Value = 0.4052844714716054
$ws.Cells.Item(10,10).Value = 0.4052844714716054
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.530935
$ws.Cells.Item(10,14).Value = 7.592805
$ws.Cells.Item(10,15).Value = 0.2497449118964618
$ws.Cells.Item(10,16).Value = 0.2497449118964618
$ws.Cells.Item(10,17).Value = 4.944790634189999
$ws.Cells.Item(10,18).Value = 44.50311570771
$ws.Cells.Item(10,19).Value = 0.1012177346206802
$ws.Cells.Item(10,20).Value = 0.1012177346206802
